# Applies the "added ifoCAST full series evaluation" update:
# Updates the error-metric table (columns B:G, rows 2-11) on the
# "first" worksheet with the new evaluation numbers.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$data = @{
    2  = @(-0.215323084493725, 2.038076249176838, 11.75110517027366, 3.427988502062639, 3.501729539135053, 22)
    3  = @(-0.3306279681019472, 1.242235062533342, 4.741407841354263, 2.177477403178794, 2.205379250257468, 21)
    4  = @(-0.1458582200366897, 0.9788677861386604, 2.521592853974958, 1.587952409228614, 1.622317456643164, 20)
    5  = @(0.05076098659158174, 0.9760318816071173, 2.328783003756465, 1.526035059805791, 1.566984368388629, 19)
    6  = @(-0.05249531032832062, 0.7762292542483971, 0.8992286170298072, 0.9482766563771393, 0.9742723218056017, 18)
    7  = @(0.06042407413349108, 0.8524080490082017, 1.247859295983782, 1.117076226577122, 1.149770081995504, 17)
    8  = @(0.1693066921153091, 0.7346782204385869, 0.7727529731107774, 0.8790636911571182, 0.8908951106930597, 16)
    9  = @(0.1421209623954074, 0.7259291398880572, 1.04014769163526, 1.019876311929667, 1.045372109252625, 15)
    10 = @(0.1894198331387246, 0.5084780843008334, 0.3620977834071827, 0.6017456135338111, 0.5927154423688374, 14)
    11 = @(0.1832594546135161, 0.7061925678054805, 0.8060402952786281, 0.8977974689642582, 0.9147828284377117, 13)
}

foreach ($row in $data.Keys) {
    $values = $data[$row]
    $ws.Cells.Item($row, 2).Value = $values[0]  # B
    $ws.Cells.Item($row, 3).Value = $values[1]  # C
    $ws.Cells.Item($row, 4).Value = $values[2]  # D
    $ws.Cells.Item($row, 5).Value = $values[3]  # E
    $ws.Cells.Item($row, 6).Value = $values[4]  # F
    $ws.Cells.Item($row, 7).Value = $values[5]  # G
}
